$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = 'Normal'
}

$ws.Range('D2').Value = '65.940.78'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '3.180.70'
$ws.Range('E3').Value = '  -0.70%  '
Set-TextValue 'D4' '1.00'
$ws.Range('E4').Value = '  -0.12%  '
Set-TextValue 'D5' '603.13'
$ws.Range('E5').Value = '  +0.85%  '
Set-TextValue 'D6' '154.02'
$ws.Range('E6').Value = '  +0.10%  '
Set-TextValue 'D7' '1.00'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.182.12'
$ws.Range('E9').Value = '  +1.88%  '
$ws.Range('E10').Value = '  -1.34%  '
Set-TextValue 'D11' '5.64'
$ws.Range('E11').Value = '  -6.93%  '
Set-TextValue 'D12' '0.512'
$ws.Range('E12').Value = '  +0.00%  '
Set-TextValue 'D13' '0.0000263'
$ws.Range('E13').Value = '  -2.95%  '
Set-TextValue 'D14' '38.09'
$ws.Range('E14').Value = '  -2.55%  '
$ws.Range('D15').Value = '3.704.16'
$ws.Range('E15').Value = '  -0.75%  '
$ws.Range('D16').Value = '66.022.44'
$ws.Range('E16').Value = '  +0.12%  '
Set-TextValue 'D17' '7.37'
$ws.Range('E17').Value = '  -0.68%  '
$ws.Range('D18').Value = '3.183.53'
$ws.Range('E18').Value = '  -0.64%  '
$ws.Range('E19').Value = '  +0.85%  '
Set-TextValue 'D20' '506.02'
$ws.Range('E20').Value = '  -0.76%  '
Set-TextValue 'D21' '15.29'
$ws.Range('E21').Value = '  -0.17%  '
Set-TextValue 'D22' '0.727'
$ws.Range('E22').Value = '  -2.01%  '
Set-TextValue 'D23' '8.00'
$ws.Range('E23').Value = '  -0.30%  '
Set-TextValue 'D24' '14.77'
$ws.Range('E24').Value = '  -3.38%  '
Set-TextValue 'D25' '84.32'
$ws.Range('E25').Value = '  -0.72%  '
$ws.Range('E26').Value = '  +0.10%  '
Set-TextValue 'D27' '2.99'
$ws.Range('E27').Value = '  -0.21%  '
Set-TextValue 'D28' '9.15'
$ws.Range('E28').Value = '  -2.15%  '
Set-TextValue 'D29' '2.37'
$ws.Range('E29').Value = '  +4.24%  '
Set-TextValue 'D30' '2.99'
$ws.Range('E30').Value = '  +4.35%  '
Set-TextValue 'D31' '7.15'
$ws.Range('E31').Value = '  +4.45%  '
Set-TextValue 'D32' '27.84'
$ws.Range('E32').Value = '  -1.68%  '
$ws.Range('E33').Value = '  +0.17%  '
Set-TextValue 'D34' '1.17'
$ws.Range('E34').Value = '  -4.54%  '
Set-TextValue 'D35' '6.46'
$ws.Range('E35').Value = '  -1.47%  '
Set-TextValue 'D36' '511.49'
$ws.Range('E36').Value = '  +5.83%  '
Set-TextValue 'D37' '55.16'
$ws.Range('E37').Value = '  +0.07%  '
Set-TextValue 'D38' '0.0879'
$ws.Range('E38').Value = '  -3.03%  '
Set-TextValue 'D39' '0.0415'
$ws.Range('E39').Value = '  -1.10%  '
$ws.Range('E40').Value = '  +4.87%  '
$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = '0.0₃0689'
$ws.Range('E41').Value = '  +6.58%  '
$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D42' '8.70'
$ws.Range('E42').Value = '  -1.47%  '
Set-TextValue 'D43' '2.84'
$ws.Range('E43').Value = '  -2.86%  '
Set-TextValue 'D44' '0.298'
$ws.Range('E44').Value = '  -0.28%  '
$ws.Range('E45').Value = '  +1.57%  '
$ws.Range('D46').Value = '2.827.57'
$ws.Range('E46').Value = '  -3.94%  '
Set-TextValue 'D47' '27.75'
$ws.Range('E47').Value = '  -2.05%  '
$ws.Range('E49').Value = '  +2.59%  '
$ws.Range('E50').Value = '  +0.20%  '
Set-TextValue 'D51' '2.60'
$ws.Range('E51').Value = '  +3.52%  '
